# Fruta / hortaliza, semanal
# The underlying edit reorders the 39 data rows (rows 2-40) of the sheet:
# every full row (columns A-R) moves to a new row position according to a
# fixed permutation, while the header row (row 1) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 40
$lastCol = 18   # column R

# Snapshot every cell's current value before moving anything, so the
# permutation can be applied safely (source rows must not be clobbered
# before they are read).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destination row -> source row (i.e. new row N gets the old content of row mapping[N])
$mapping = @{
    2  = 34
    3  = 7
    4  = 12
    5  = 13
    6  = 17
    7  = 30
    8  = 29
    9  = 32
    10 = 40
    11 = 38
    12 = 31
    13 = 2
    14 = 25
    15 = 11
    16 = 20
    17 = 37
    18 = 35
    19 = 36
    20 = 33
    21 = 14
    22 = 21
    23 = 5
    24 = 26
    25 = 24
    26 = 28
    27 = 39
    28 = 15
    29 = 10
    30 = 4
    31 = 22
    32 = 27
    33 = 6
    34 = 23
    35 = 18
    36 = 16
    37 = 9
    38 = 19
    39 = 8
    40 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
